$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the header text in B1 from "descripciones" to "descripcion"
$ws.Range("B1").Value = "descripcion"

# Move / update the active selection to B2, as left by the author after editing B1
$ws.Range("B2").Select() | Out-Null
